$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Insert the two new rows first (before touching hyperlinks),
#    so hyperlink refs we add later land on their final row.
# ---------------------------------------------------------------
# New row 5 "DIO_Layered Arch." goes in right after the current row 4 (DIO)
$ws.Rows(5).Insert()
# New row 7 "keypad" goes in right after the current row 6 (LCD, after the
# previous insert shifted it down from 5).
$ws.Rows(7).Insert()

# ---------------------------------------------------------------
# 2) Update existing cell text that changed.
# ---------------------------------------------------------------
$ws.Range("B4").Value = "relay,buzzer,Alarm app"

# ---------------------------------------------------------------
# 3) Fill in the two brand-new rows (A & B columns).
# ---------------------------------------------------------------
$ws.Range("A5").Value = "DIO_Layered Arch."
$ws.Range("B5").Value = "DIO Driver with Layered Architecture"

$ws.Range("A7").Value = "keypad"
$ws.Range("B7").Value = "KP layered Arch, kP with LCD app"

# ---------------------------------------------------------------
# 4) Column C - widen it and add the Material hyperlinks.
# ---------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 103.16666666666667

# C2 / C3 need a display text that differs from the actual cell text
# (spaces instead of %20), so add the hyperlink first (which seeds both
# the cell text and the display with the spaced text) and then overwrite
# the cell's Value afterwards - that keeps the hyperlink's display intact
# while fixing the underlying cell text to the encoded URL.
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/computer%20architecture%20ver%201.3.pdf", "", "", "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/computer architecture ver 1.3.pdf")
$ws.Range("C2").Value = "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/computer%20architecture%20ver%201.3.pdf"

$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/embedded%20C.pptx", "", "", "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/embedded C.pptx")
$ws.Range("C3").Value = "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/embedded%20C.pptx"

$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/DIO_Interfacing.pptx")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/GPIO.pptx")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/LCD.pptx")

# Interrupt row (row 8) got its hyperlink before the keypad row's.
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/Interrupt.pptx")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/Keypad.pptx")

$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/ADC.pptx")

# C10 (Timer) gets the hyperlink-look text + style but no live hyperlink:
# add it and then delete just that one link, leaving the text/style behind.
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/Timer.pptx")
$ws.Range("C10").Hyperlinks.Item(1).Delete()

$ws.Hyperlinks.Add($ws.Range("C12"), "https://github.com/Mohamedsaied8/AMIT_labs/blob/master/slides/UART.pptx")

# ---------------------------------------------------------------
# 5) Final selection, matching where editing ended up.
# ---------------------------------------------------------------
$ws.Range("C16").Select() | Out-Null
